# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Arveja Verde" (Vega Modelo de Temuco)
# above the existing row 57, shifting the old rows 57-59 down to 59-61.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 57 (pushes current rows 57..59 down to 59..61)
$ws.Rows.Item(57).Insert()
$ws.Rows.Item(57).Insert()

# --- New row 57 ---
$ws.Cells.Item(57,1).Value2 = 10
$ws.Cells.Item(57,2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(57,3).Value2 = "La Araucanía"
$ws.Cells.Item(57,4).Value2 = 44516
$ws.Cells.Item(57,5).Value2 = 9
$ws.Cells.Item(57,6).Value2 = 100112022
$ws.Cells.Item(57,7).Value2 = "Arveja Verde"
$ws.Cells.Item(57,8).Value2 = "Sin especificar"
$ws.Cells.Item(57,9).Value2 = "Primera"
$ws.Cells.Item(57,10).Value2 = 35
$ws.Cells.Item(57,11).Value2 = 15000
$ws.Cells.Item(57,12).Value2 = 15000
$ws.Cells.Item(57,13).Value2 = 15000
$ws.Cells.Item(57,14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(57,15).Value2 = "Región Metropolitana"
$ws.Cells.Item(57,16).Value2 = 600
$ws.Cells.Item(57,17).Value2 = 25
$ws.Cells.Item(57,18).Value2 = "Hortaliza"

# --- New row 58 ---
$ws.Cells.Item(58,1).Value2 = 10
$ws.Cells.Item(58,2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(58,3).Value2 = "La Araucanía"
$ws.Cells.Item(58,4).Value2 = 44516
$ws.Cells.Item(58,5).Value2 = 9
$ws.Cells.Item(58,6).Value2 = 100112022
$ws.Cells.Item(58,7).Value2 = "Arveja Verde"
$ws.Cells.Item(58,8).Value2 = "Sin especificar"
$ws.Cells.Item(58,9).Value2 = "Primera"
$ws.Cells.Item(58,10).Value2 = 55
$ws.Cells.Item(58,11).Value2 = 16000
$ws.Cells.Item(58,12).Value2 = 16000
$ws.Cells.Item(58,13).Value2 = 16000
$ws.Cells.Item(58,14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(58,15).Value2 = "Región del Maule"
$ws.Cells.Item(58,16).Value2 = 640
$ws.Cells.Item(58,17).Value2 = 25
$ws.Cells.Item(58,18).Value2 = "Hortaliza"
